$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 257, shifting the existing rows 257:270 down to 258:271
$ws.Rows("257:257").Insert()

# Populate the new row 257 with the new weekly price observation
$ws.Cells.Item(257, 1).Value = 8
$ws.Cells.Item(257, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(257, 3).Value = "Coquimbo"
$ws.Cells.Item(257, 4).Value = "2023-12-07"
$ws.Cells.Item(257, 5).Value = 4
$ws.Cells.Item(257, 6).Value = 100112044
$ws.Cells.Item(257, 7).Value = "Perejil"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 2000
$ws.Cells.Item(257, 11).Value = 2300
$ws.Cells.Item(257, 12).Value = 2500
$ws.Cells.Item(257, 13).Value = 2400
$ws.Cells.Item(257, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(257, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(257, 16).Value = 1600
$ws.Cells.Item(257, 17).Value = 1.5
$ws.Cells.Item(257, 18).Value = "Hortaliza"
